$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9796509146690369
$ws.Range("B1").Value = 0.8403860330581665
$ws.Range("C1").Value = 2.443525552749634
$ws.Range("D1").Value = 5.228856086730957
$ws.Range("E1").Value = 1.152986645698547
